$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2751.5881
$ws.Range("J40").Value = 3340.2666
$ws.Range("L40").Value = 3340.2666
$ws.Range("N40").Value = -3690.2666

$ws.Range("H51").Value = 12618.091
$ws.Range("I51").Value = 27650
$ws.Range("J51").Value = 4028.4285
$ws.Range("K51").Value = 27650
$ws.Range("L51").Value = 4028.4285
$ws.Range("M51").Value = -27166
$ws.Range("N51").Value = -4996.4285

$ws.Range("H98").Value = 2217.6
$ws.Range("I98").Value = 2217.6
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2217.6
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -719.5999999999999
$ws.Range("N98").ClearContents()

$ws.Range("H116").Value = 2580
$ws.Range("J116").Value = 2500
$ws.Range("L116").Value = 2500
$ws.Range("N116").Value = -9384

$ws.Range("H122").Value = 2217.6
$ws.Range("I122").Value = 2217.6
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6652.799999999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4202.799999999999
$ws.Range("N122").ClearContents()

$ws.Range("H138").Value = 3850.68
$ws.Range("I138").Value = 2939.1
$ws.Range("J138").Value = 4458.4
$ws.Range("K138").Value = 8817.299999999999
$ws.Range("L138").Value = 13375.2
$ws.Range("M138").Value = -3677.299999999999
$ws.Range("N138").Value = -23655.2

$ws.Range("H141").Value = 4670.7144
$ws.Range("I141").Value = 4723.75
$ws.Range("J141").Value = 4600
$ws.Range("K141").Value = 14171.25
$ws.Range("L141").Value = 13800
$ws.Range("M141").Value = -8991.25
$ws.Range("N141").Value = -24160

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H103").Value = 42356
$ws.Range("J103").Value = 42356
$ws.Range("L103").Value = 42356
$ws.Range("N103").Value = -44700

$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 297.5
$ws.Range("I22").Value = 297.5
$ws.Range("K22").Value = 297.5
$ws.Range("M22").Value = -124.5

$ws.Range("H140").Value = 67333
$ws.Range("J140").Value = 67333
$ws.Range("L140").Value = 67333
$ws.Range("N140").Value = -77693

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21849.986
$ws.Range("I31").Value = 1405.0385
$ws.Range("J31").Value = 33662.62
$ws.Range("K31").Value = 1405.0385
$ws.Range("L31").Value = 33662.62
$ws.Range("M31").Value = -1110.0385
$ws.Range("N31").Value = -34252.62

$ws.Range("H34").Value = 21849.986
$ws.Range("I34").Value = 1405.0385
$ws.Range("J34").Value = 33662.62
$ws.Range("K34").Value = 1405.0385
$ws.Range("L34").Value = 33662.62
$ws.Range("M34").Value = -1203.0385
$ws.Range("N34").Value = -34066.62

$ws.Range("H86").Value = 3880.8
$ws.Range("J86").Value = 4184.6665
$ws.Range("L86").Value = 4184.6665
$ws.Range("N86").Value = -6430.6665

$ws.Range("H89").Value = 3880.8
$ws.Range("J89").Value = 4184.6665
$ws.Range("L89").Value = 20923.3325
$ws.Range("N89").Value = -32155.3325

$ws.Range("H106").Value = 32499.5
$ws.Range("J106").Value = 32499.5
$ws.Range("L106").Value = 32499.5
$ws.Range("N106").Value = -35023.5

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 50
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 300
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -526

$ws.Range("H5").Value = 1114.6571
$ws.Range("I5").Value = 788.9375
$ws.Range("K5").Value = 2366.8125
$ws.Range("M5").Value = -2254.8125

$ws.Range("H37").Value = 930776.4399999999
$ws.Range("J37").Value = 930776.4399999999
$ws.Range("L37").Value = 2792329.32
$ws.Range("N37").Value = -2792553.32

$ws.Range("H38").Value = 37.23077
$ws.Range("I38").Value = 23.75
$ws.Range("K38").Value = 71.25
$ws.Range("M38").Value = 275.75

$ws.Range("H98").Value = 200475.8
$ws.Range("I98").Value = 595
$ws.Range("K98").Value = 1785
$ws.Range("M98").Value = -287

$ws.Range("H131").Value = 780.99
$ws.Range("I131").Value = 338.66666
$ws.Range("J131").Value = 809.2234
$ws.Range("K131").Value = 1015.99998
$ws.Range("L131").Value = 2427.6702
$ws.Range("M131").Value = 4024.00002
$ws.Range("N131").Value = -12507.6702

$ws.Range("H133").Value = 4050
$ws.Range("I133").Value = 3100
$ws.Range("J133").Value = 5000
$ws.Range("K133").Value = 9300
$ws.Range("L133").Value = 15000
$ws.Range("M133").Value = -4240
$ws.Range("N133").Value = -25120

$ws.Range("H135").Value = 1114.6571
$ws.Range("I135").Value = 788.9375
$ws.Range("K135").Value = 7100.4375
$ws.Range("M135").Value = -4565.4375

$ws.Range("H138").Value = 11266.182
$ws.Range("I138").Value = 14366.125
$ws.Range("K138").Value = 43098.375
$ws.Range("M138").Value = -37958.375

$ws.Range("H139").Value = 2250.08
$ws.Range("I139").Value = 1256
$ws.Range("J139").Value = 2912.8
$ws.Range("K139").Value = 3768
$ws.Range("L139").Value = 8738.400000000001
$ws.Range("M139").Value = 1372
$ws.Range("N139").Value = -19018.4

$ws.Range("H140").Value = 5946.75
$ws.Range("I140").Value = 9152.691999999999
$ws.Range("K140").Value = 27458.076
$ws.Range("M140").Value = -22278.076

$ws.Range("H141").Value = 7913.75
$ws.Range("I141").Value = 9064.166999999999
$ws.Range("J141").Value = 4462.5
$ws.Range("K141").Value = 27192.501
$ws.Range("L141").Value = 13387.5
$ws.Range("M141").Value = -22012.501
$ws.Range("N141").Value = -23747.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 9068078
$ws.Range("J16").Value = 1669163.9
$ws.Range("L16").Value = 1669163.9
$ws.Range("N16").Value = -1669503.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 122772
$ws.Range("I81").Value = 691.5
$ws.Range("K81").Value = 1383
$ws.Range("M81").Value = -322

$ws.Range("H84").Value = 122772
$ws.Range("I84").Value = 691.5
$ws.Range("K84").Value = 6915
$ws.Range("M84").Value = -1611
